$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.671.03"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").Value = "1.697.80"
$ws.Range("E3").Value = "  +0.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'316.17"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.18%  "

# Row 7
$ws.Range("D7").Value = "'0.3926"
$ws.Range("E7").Value = "  -0.39%  "

# Row 8
$ws.Range("D8").Value = "'0.4051"
$ws.Range("E8").Value = "  +0.74%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").Value = "'1.003"
$ws.Range("E10").Value = "  +0.20%  "

# Row 11
$ws.Range("D11").Value = "'52.99"
$ws.Range("E11").Value = "  -1.47%  "

# Row 12
$ws.Range("D12").Value = "'0.08845"
$ws.Range("E12").Value = "  +0.83%  "

# Row 13
$ws.Range("D13").Value = "'7.426"
$ws.Range("E13").Value = "  +2.70%  "

# Row 14
$ws.Range("D14").Value = "'23.66"
$ws.Range("E14").Value = "  +1.88%  "

# Row 15
$ws.Range("D15").Value = "'8.108"
$ws.Range("E15").Value = "  +6.88%  "

# Row 16
$ws.Range("D16").Value = "'0.00001319"
$ws.Range("E16").Value = "  -0.52%  "

# Row 17
$ws.Range("D17").Value = "1.699.67"
$ws.Range("E17").Value = "  +0.12%  "

# Row 18
$ws.Range("D18").Value = "'99.37"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19
$ws.Range("D19").Value = "'0.07035"
$ws.Range("E19").Value = "  -0.49%  "

# Row 21
$ws.Range("D21").Value = "'7.071"
$ws.Range("E21").Value = "  +2.92%  "

# Row 22
$ws.Range("E22").Value = "  +0.49%  "

# Row 23
$ws.Range("D23").Value = "'14.72"
$ws.Range("E23").Value = "  +4.87%  "

# Row 24
$ws.Range("D24").Value = "24.657.32"
$ws.Range("E24").Value = "  +0.28%  "

# Row 25
$ws.Range("D25").Value = "'3.134"
$ws.Range("E25").Value = "  +3.00%  "

# Row 26
$ws.Range("D26").Value = "'2.349"
$ws.Range("E26").Value = "  +1.57%  "

# Row 27
$ws.Range("D27").Value = "'22.62"
$ws.Range("E27").Value = "  +1.02%  "

# Row 28
$ws.Range("D28").Value = "'163.68"
$ws.Range("E28").Value = "  +2.30%  "

# Row 29
$ws.Range("D29").Value = "'8.781"
$ws.Range("E29").Value = "  +18.00%  "

# Row 30
$ws.Range("D30").Value = "'135.60"

# Row 31
$ws.Range("D31").Value = "'5.142"
$ws.Range("E31").Value = "  -1.72%  "

# Row 32
$ws.Range("D32").Value = "'0.09004"
$ws.Range("E32").Value = "  +5.55%  "

# Row 33
$ws.Range("D33").Value = "'7.601"
$ws.Range("E33").Value = "  +4.89%  "

# Row 34
$ws.Range("D34").Value = "'1.070"
$ws.Range("E34").Value = "  -3.24%  "

# Row 35
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02961"
$ws.Range("E35").Value = "  +7.47%  "

# Row 36
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.967"
$ws.Range("E36").Value = "  +0.69%  "

# Row 37
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'11.04"
$ws.Range("E37").Value = "  -2.74%  "

# Row 38
$ws.Range("D38").Value = "'0.2752"
$ws.Range("E38").Value = "  +0.70%  "

# Row 39
$ws.Range("E39").Value = "  -0.49%  "

# Row 40
$ws.Range("D40").Value = "'0.09173"
$ws.Range("E40").Value = "  +1.26%  "

# Row 41
$ws.Range("D41").Value = "'1.458"
$ws.Range("E41").Value = "  -0.31%  "

# Row 42
$ws.Range("D42").Value = "'0.7665"
$ws.Range("E42").Value = "  -0.72%  "

# Row 43
$ws.Range("D43").Value = "'16.02"
$ws.Range("E43").Value = "  +3.66%  "

# Row 44
$ws.Range("D44").Value = "'2.588"
$ws.Range("E44").Value = "  +1.92%  "

# Row 45
$ws.Range("D45").Value = "'0.7172"
$ws.Range("E45").Value = "  -0.48%  "

# Row 46
$ws.Range("D46").Value = "'4.214"
$ws.Range("E46").Value = "  +0.09%  "

# Row 47
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48
$ws.Range("D48").Value = "'1.339"
$ws.Range("E48").Value = "  -1.16%  "

# Row 49
$ws.Range("D49").Value = "'139.83"
$ws.Range("E49").Value = "  -1.08%  "

# Row 50
$ws.Range("D50").Value = "'0.07976"
$ws.Range("E50").Value = "  -0.57%  "

# Row 51
$ws.Range("D51").Value = "'90.24"
$ws.Range("E51").Value = "  +1.98%  "
